# Weekly refresh of the "Poroto verde" sheet.
# The underlying data rows (identified by their Fecha/Precio/Origen block)
# get reshuffled to their new weekly positions while the rest of each row
# (Mercado, Region, Codreg, Categoria, Calidad, Kg/Unidades, Clasificacion)
# stays put. Columns touched per row: D, H, J, K, L, M, N, O, P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 32
$cols = @("D", "H", "J", "K", "L", "M", "N", "O", "P")

# Row r (destination) receives the block of values that currently (before
# this edit) lives in row $rowMap[r] (source).
$rowMap = @{
    2 = 30; 3 = 19; 4 = 27; 5 = 4;  6 = 8;  7 = 26; 8 = 23; 9 = 25; 10 = 31;
    11 = 17; 12 = 15; 13 = 2; 14 = 22; 15 = 7; 16 = 20; 17 = 12; 18 = 13;
    19 = 3; 20 = 9; 21 = 5; 22 = 24; 23 = 32; 24 = 29; 25 = 28; 26 = 6;
    27 = 16; 28 = 10; 29 = 18; 30 = 21; 31 = 11; 32 = 14
}

# 1) Snapshot the current values of the affected columns for every row
#    before any writes happen, so later writes don't clobber data we
#    still need to read for other rows.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# 2) Write back each destination row using the snapshotted source row data.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $rowMap[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}
